# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.481.90'
$ws.Range("D3").Value = '1.919.12'
$ws.Range("E3").Value = '  +1.65%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.19'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4740'
$ws.Range("E7").Value = '  +2.49%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4095'
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.84'
$ws.Range("E9").Value = '  +0.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08047'
$ws.Range("E10").Value = '  +0.76%  '
$ws.Range("E11").Value = '  +2.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.57'
$ws.Range("E12").Value = '  +3.91%  '
$ws.Range("D13").Value = '1.930.11'
$ws.Range("E13").Value = '  +2.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.938'
$ws.Range("E14").Value = '  +0.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.164'
$ws.Range("E15").Value = '  +1.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '89.78'
$ws.Range("E16").Value = '  +0.74%  '
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06594'
$ws.Range("E18").Value = '  +0.43%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001031'
$ws.Range("E19").Value = '  +0.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.77'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("D22").Value = '29.499.52'
$ws.Range("E22").Value = '  +1.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.547'
$ws.Range("E23").Value = '  +3.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.52'
$ws.Range("E24").Value = '  +2.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.208'
$ws.Range("E25").Value = '  -0.27%  '
$ws.Range("D26").Value = '2.150.44'
$ws.Range("E26").Value = '  +2.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '155.23'
$ws.Range("E27").Value = '  -1.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.85'
$ws.Range("E28").Value = '  +0.95%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.987'
$ws.Range("E29").Value = '  +10.66%  '
$ws.Range("E30").Value = '  +0.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '117.92'
$ws.Range("E31").Value = '  -0.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.067'
$ws.Range("E32").Value = '  +9.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09538'
$ws.Range("E33").Value = '  +2.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.437'
$ws.Range("E34").Value = '  +1.47%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.558'
$ws.Range("E35").Value = '  -1.20%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.416'
$ws.Range("E36").Value = '  +2.62%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06115'
$ws.Range("E37").Value = '  +0.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02261'
$ws.Range("E38").Value = '  +1.29%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.328'
$ws.Range("E39").Value = '  +0.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.173'
$ws.Range("E40").Value = '  -0.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5911'
$ws.Range("E41").Value = '  +2.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.561'
$ws.Range("E42").Value = '  +11.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1848'
$ws.Range("E43").Value = '  +1.60%  '
$ws.Range("E44").Value = '  +0.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.08016'
$ws.Range("E45").Value = '  +14.34%  '
$ws.Range("E46").Value = '  +1.92%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5568'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '12.18'
$ws.Range("E48").Value = '  +0.93%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.936'
$ws.Range("E49").Value = '  +1.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '113.11'
$ws.Range("E50").Value = '  +2.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.85'
$ws.Range("E51").Value = '  +0.97%  '
